# Fruta / hortaliza, semanal
#
# A new weekly price record is inserted at the top of the Berenjena /
# Macroferia Regional de Talca data block (row 146). All the existing
# records that used to live in rows 146-160 are pushed down by one row
# (146->147, 147->148, ..., 159->160, 160->161), and the now-empty row 146
# is filled in with the new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 146
$lastRow = 160

# Shift existing rows down by one, starting from the bottom so that a row
# is never overwritten before it has been copied to its new location.
for ($r = $lastRow; $r -ge $firstRow; $r--) {
    $dstRow = $r + 1
    $srcRange = $ws.Range("A$r`:R$r")
    $dstRange = $ws.Range("A$dstRow`:R$dstRow")

    $dstRange.Value2 = $srcRange.Value2

    # Preserve the date number format used on column D.
    $ws.Range("D$dstRow").NumberFormat = $ws.Range("D$r").NumberFormat
}

# Fill row 146 in with the new record.
$ws.Range("D146").Value2 = 45013
$ws.Range("K146").Value2 = 8000
$ws.Range("L146").Value2 = 8000
$ws.Range("M146").Value2 = 8000
$ws.Range("P146").Value2 = 160
